$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows before the "Total (hrs)" summary rows (old rows 22 & 23
# shift down to become rows 25 & 26).
$ws.Rows("22:24").Insert()

# --- Row 22: new "Development" log entry (2019-06-09) ---
$ws.Range("B22").Value = "Development"
$ws.Range("C22").Value = 43625
$ws.Range("D22").Value = 10
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = "1) Basic implementation of audio webworklet for FFT processing with window size of 512 samples."
$ws.Rows("22:22").RowHeight = 29

# --- Row 23: new "Development" log entry (2019-06-10) ---
$ws.Range("B23").Value = "Development"
$ws.Range("C23").Value = 43626
$ws.Range("D23").Value = 6
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = "1) Converted use of Promises to async/await for code clarity in webworklet & subtitles.js`n2) Researched on TextTracks cues system`n3) Fixed subtitles bug involving Panopto's desynced two videos system (videos play at different offsets at different timings, only synced at runtime in implementation. Our implementation must thus also be synced at runtime)"
$ws.Rows("23:23").RowHeight = 87

# --- Row 24: new "Development" log entry (2019-06-11) ---
$ws.Range("B24").Value = "Development"
$ws.Range("C24").Value = 43627
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = "1) Implemented silence detection using TextTracks oncuechange (cues system)`n2) Brainstorming for seeking videos (need to keep videos in sync)`n3) Discovered bug involving lag when seeking, bug appears to be involving buffering with HLSJS`n4) More research into Panopto's implementation, HLS.js documentation on their API, Debugging`n5) Implemented basic implementation of silence detection"
$ws.Rows("24:24").RowHeight = 116

# --- Row 26 (previously row 23): update summary formulas to cover the new rows ---
$ws.Range("C26").Formula = "=SUM(D26:E26)"
$ws.Range("D26").Formula = "=SUM(D3:D24)"
$ws.Range("E26").Formula = "=SUM(E3:E24)"

# --- Update the view state to match where the author last left the sheet ---
$ws.Range("A24").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F25").Select() | Out-Null
